# 2des pwbe aulas12 e 13
# Adds attendance column CB (aula 13, 2022-11-04) values alongside the
# existing CA column (aula 12), for every student row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chamada")

$attendance = [ordered]@{
    3  = "P"
    4  = "P"
    5  = "F"
    6  = "P"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    12 = "P"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "P"
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "F"
    29 = "P"
    30 = "P"
}

foreach ($row in $attendance.Keys) {
    $cell = $ws.Range("CB$row")
    $cell.HorizontalAlignment = -4131
    $cell.Value = $attendance[$row]
}

$ws.Range("CB26").Select()
